$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22: Team 2 changes from "Maharashtra Sports Club" to "Westridge Warriors"
$ws.Range("D22").Value = "Westridge Warriors"

# Row 22: Venue cell is cleared
$ws.Range("E22").ClearContents()

# Text no longer wraps to two lines, so the row's explicit height is dropped
# back to the sheet's automatic/default height.
$ws.Rows.Item(22).EntireRow.AutoFit()

# Update selection to E22 to match the saved cursor position
$ws.Range("E22").Select()
